$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores its figures as literal text in the source data
# (e.g. "599.44"), not as numbers. Force Text format before writing so Excel
# does not auto-convert a numeric-looking string to a number -- which would both
# change the cell type and, for a case like "5.00", silently drop the trailing
# zero (-> 5). The "Volume(1h)" column (E) always contains a "%" and spaces, so
# it is never at risk of numeric coercion and needs no special handling.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.935.31"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.553.41"
$ws.Range("E3").Value = "  +4.60%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.44"
$ws.Range("E5").Value = "  +3.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.75"
$ws.Range("E6").Value = "  +3.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.552.05"
$ws.Range("E7").Value = "  +4.55%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +3.72%  "

$ws.Range("E10").Value = "  +3.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.93"
$ws.Range("E11").Value = "  -0.12%  "

$ws.Range("E12").Value = "  +4.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.154.51"
$ws.Range("E13").Value = "  +4.55%  "

$ws.Range("E14").Value = "  +4.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.14"
$ws.Range("E15").Value = "  +5.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.555.46"
$ws.Range("E16").Value = "  +5.55%  "

$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.857.83"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.04"
$ws.Range("E19").Value = "  +8.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.42"
$ws.Range("E20").Value = "  +8.06%  "

$ws.Range("E21").Value = "  +3.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "389.26"
$ws.Range("E22").Value = "  +3.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.576"
$ws.Range("E23").Value = "  +7.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.699.24"
$ws.Range("E24").Value = "  +4.78%  "

$ws.Range("E25").Value = "  +3.85%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").Value = "  +14.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.66"
$ws.Range("E28").Value = "  +8.41%  "

$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("E30").Value = "  +5.80%  "

$ws.Range("E31").Value = "  +6.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.562.35"
$ws.Range("E32").Value = "  +4.26%  "

$ws.Range("E33").Value = "  +22.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.01"
$ws.Range("E34").Value = "  +5.46%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("E36").Value = "  +3.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "170.32"
$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.95"
$ws.Range("E38").Value = "  +5.85%  "

$ws.Range("E39").Value = "  +8.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.00"
$ws.Range("E40").Value = "  +10.13%  "

$ws.Range("E41").Value = "  +7.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.828"
$ws.Range("E42").Value = "  +4.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.68"
$ws.Range("E43").Value = "  +21.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.64"
$ws.Range("E44").Value = "  +2.29%  "

$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("E46").Value = "  +6.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("E47").Value = "  +11.16%  "

$ws.Range("E48").Value = "  +4.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.446.57"
$ws.Range("E49").Value = "  +12.36%  "

$ws.Range("E50").Value = "  +7.27%  "

$ws.Range("E51").Value = "  +18.43%  "
